$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the wind-speed (U [m/s]) values in column B, rows 13-38
# Rows 13-23 (Time 11-21 s): 15 -> 10
$ws.Range("B13:B23").Value = 10

# Rows 24-27 (Time 22-25 s): 15 -> 0
$ws.Range("B24:B27").Value = 0

# Rows 28-32 (Time 26-30 s): 6 -> 0
$ws.Range("B28:B32").Value = 0

# Rows 33-38 (Time 31-36 s): 6 -> 15
$ws.Range("B33:B38").Value = 15

# Scroll the view down and move the selection, matching the saved view state
$excel.ActiveWindow.ScrollRow = 11
$ws.Range("F37").Select()
